# Update the "Förändrad" (changed) date column C for rows 2-24
# from serial date 45229 (2023-10-30) to 45231 (2023-11-01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45229) {
        $cell.Value = 45231
    }
}
